$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Add the two missing "Slides" links for class_3 and class_4 into column G
# (same text already used for Sheet2!B4 and Sheet2!B5).
$ws1.Range("G4").Value = "[Slides](slides/class_3/class_3#1) [.qmd](slides/class_3/class_3.qmd)"
$ws1.Range("G5").Value = "[Slides](slides/class_4/class_4#1) [.qmd](slides/class_4/class_4.qmd) [.R](slides/class_4/class_4_taller.R)"

# Update the active selections on both sheets, as recorded in the sheetViews.
# Sheet2 is updated first, then Sheet1 last so Sheet1 stays the active
# (tab-selected) sheet, matching the workbook's original state.
$ws2.Activate() | Out-Null
$ws2.Range("B4:B5").Select() | Out-Null
$ws2.Cells.Item(5, 2).Activate() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("G4:G5").Select() | Out-Null
